# Add a "Status" column (F) with "APROVADO" to both test-case sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("C.T - Válido")
$ws2 = $wb.Worksheets.Item("C.T - Inválido")

# ---------------------------------------------------------------------
# Sheet 1: "C.T - Válido"  (rows 2-7, header + 5 data rows)
# ---------------------------------------------------------------------

# Header cell F2 - copy formatting from the existing "Status"-like header
# style used on sheet2!E2 (blue fill header, no top/bottom border) so we
# reuse the existing style record instead of minting a new one.
$ws2.Range("E2").Copy()
$ws1.Range("F2").PasteSpecial(-4122)
$ws1.Range("F2").Value = "Status"

# Data cells F3:F7 - copy formatting from an existing plain data cell
# (B3, style used for normal body cells) then write the value.
$ws1.Range("B3").Copy()
$ws1.Range("F3:F7").PasteSpecial(-4122)
$ws1.Range("F3").Value = "APROVADO"
$ws1.Range("F4").Value = "APROVADO"
$ws1.Range("F5").Value = "APROVADO"
$ws1.Range("F6").Value = "APROVADO"
$ws1.Range("F7").Value = "APROVADO"

# New column F width (~10 chars, matches the authored width)
$ws1.Columns("F").ColumnWidth = 9.14

# ---------------------------------------------------------------------
# Sheet 2: "C.T - Inválido"  (rows 2-6, header + 4 data rows)
# ---------------------------------------------------------------------

# Header cell F2 - same header style (already present at E2 on this sheet).
$ws2.Range("E2").Copy()
$ws2.Range("F2").PasteSpecial(-4122)
$ws2.Range("F2").Value = "Status"

# Data cells F3:F6
$ws1.Range("B3").Copy()
$ws2.Range("F3:F6").PasteSpecial(-4122)
$ws2.Range("F3").Value = "APROVADO"
$ws2.Range("F4").Value = "APROVADO"
$ws2.Range("F5").Value = "APROVADO"
$ws2.Range("F6").Value = "APROVADO"

# New column F width (slightly wider on this sheet per the authored file)
$ws2.Columns("F").ColumnWidth = 10

# ---------------------------------------------------------------------
# Selections / active sheet.
# Select on the sheet that should stay inactive first, then finish on
# the sheet/cell that should end up active - last Select() wins for
# both the sheet-level tabSelected flag and the workbook's activeTab.
# ---------------------------------------------------------------------

$ws2.Range("F2:F6").Select()
$ws1.Range("C16").Select()
